# Insert a new data row at row 314 (shifts existing rows 314-417 down to 315-418)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(314).Insert()

$ws.Range("A314").Value = 10
$ws.Range("B314").Value = "Vega Modelo de Temuco"
$ws.Range("C314").Value = "La Araucanía"
$ws.Range("D314").Value = 44588
$ws.Range("D314").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E314").Value = 9
$ws.Range("F314").Value = 100112043
$ws.Range("G314").Value = "Pepino ensalada"
$ws.Range("H314").Value = "Sin especificar"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 280
$ws.Range("K314").Value = 11000
$ws.Range("L314").Value = 12000
$ws.Range("M314").Value = 11446
$ws.Range("N314").Value = "`$/caja 60 unidades"
$ws.Range("O314").Value = "Región del Maule"
$ws.Range("P314").Value = 191
$ws.Range("Q314").Value = 60
$ws.Range("R314").Value = "Hortaliza"
